# Update "paises" (countries) COVID data sheet:
#  1) Refresh the case counters for the countries whose figures changed.
#  2) Re-sort the full country table by "Casos totales" (column B) descending,
#     since ranking shifts for a few countries change their row order.
#  3) Bump the "last updated" timestamp footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 4
$lastDataRow = 219

# New figures: Country -> Casos totales, Nuevos casos, Casos activos, Recuperados,
#              Casos criticos, Muertes hoy, Muertes
$updates = @{
    "Estados Unidos"      = @(4203872, 33554, 1991969, 2064067, 0, 487, 147836)
    "India"               = @(1337021, 48891, 848986, 456630, 0, 760, 31405)
    "Chile"               = @(341304, 2545, 313696, 18770, 0, 0, 8838)
    "España"              = @(319501, 2255, 0, 0, 0, 3, 28432)
    "Alemania"            = @(205551, 409, 189400, 6959, 0, 5, 9192)
    "Francia"             = @(180528, 1130, 80472, 69864, 0, 10, 30192)
    "Irlanda"             = @(25845, 19, 23364, 718, 0, 0, 1763)
    "Marruecos"           = @(18834, 570, 16100, 2435, 0, 7, 299)
    "Zambia"              = @(3856, 67, 1677, 2043, 0, 2, 136)
    "Libano"              = @(3407, 147, 1666, 1695, 0, 3, 46)
    "Maldivas"            = @(3175, 55, 2498, 662, 0, 0, 15)
    "Cabo Verde"          = @(2220, 30, 1216, 982, 0, 1, 22)
    "Suazilandia"         = @(2073, 52, 929, 1116, 0, 0, 28)
    "Yemen"               = @(1674, 20, 779, 426, 0, 8, 469)
    "Jordania"            = @(1146, 15, 1035, 100, 0, 0, 11)
    "Burundi"             = @(345, 0, 270, 74, 0, 0, 10)
    "Birmania"            = @(346, 3, 286, 54, 0, 0, 6)
    "Monaco"              = @(116, 2, 100, 12, 0, 0, 4)
    "Papua Nueva Guinea"  = @(32, 1, 11, 21, 0, 0, 0)
    "Curazao"             = @(29, 1, 24, 4, 0, 0, 1)
}

# Build a lookup of country name -> row number for the current (pre-sort) layout.
$rowByCountry = @{}
for ($i = $firstDataRow; $i -le $lastDataRow; $i++) {
    $name = $ws.Cells.Item($i, 1).Value()
    $rowByCountry[$name] = $i
}

foreach ($country in $updates.Keys) {
    $r = $rowByCountry[$country]
    if ($r -eq $null) {
        continue
    }
    $vals = $updates[$country]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
    $ws.Cells.Item($r, 7).Value = $vals[5]
    $ws.Cells.Item($r, 8).Value = $vals[6]
}

# Re-sort the whole data block by "Casos totales" (column B), descending.
$dataRange = $ws.Range("A" + $firstDataRow + ":H" + $lastDataRow)
$keyRange = $ws.Range("B" + $firstDataRow + ":B" + $lastDataRow)
$dataRange.Sort($keyRange, 2)

# "Groenlandia" and "Islas Malvinas" are tied on every numeric column, so a
# stable sort alone leaves them in their original relative order. The source
# data swapped their relative order for this refresh, so fix that pair up
# explicitly once they're adjacent to each other post-sort.
$rowByCountry2 = @{}
for ($i = $firstDataRow; $i -le $lastDataRow; $i++) {
    $name = $ws.Cells.Item($i, 1).Value()
    $rowByCountry2[$name] = $i
}
$rGroen = $rowByCountry2["Groenlandia"]
$rMalv = $rowByCountry2["Islas Malvinas"]
if (($rGroen -ne $null) -and ($rMalv -ne $null) -and ($rGroen -gt $rMalv)) {
    for ($c = 1; $c -le 8; $c++) {
        $tmp = $ws.Cells.Item($rMalv, $c).Value()
        $ws.Cells.Item($rMalv, $c).Value = $ws.Cells.Item($rGroen, $c).Value()
        $ws.Cells.Item($rGroen, $c).Value = $tmp
    }
}

# Update the "last updated" footer text.
$ws.Range("A1").Value = "Datos actualizados a 24 de Julio de 2020 a las 20:10"
